$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the "Csakvari vs Siofok" / "Szeged vs Budafoki" match data between
# --- row 124 and row 126 (only the match-specific columns change; A-E, G, I
# --- stay where they are since the Indice/date fields are unaffected).

# New values for row 124 (previously held by row 126)
$ws.Cells.Item(124, 6).Value = "Szeged"
$ws.Cells.Item(124, 8).Value = "Budafoki"
$ws.Cells.Item(124, 10).Value = 1.86
$ws.Cells.Item(124, 12).Value = 1.91
$ws.Cells.Item(124, 13).Value = "05/11/2023 16:44"
$ws.Cells.Item(124, 14).Value = 3.25
$ws.Cells.Item(124, 16).Value = 3.31
$ws.Cells.Item(124, 17).Value = "05/11/2023 16:44"
$ws.Cells.Item(124, 18).Value = 3.88
$ws.Cells.Item(124, 20).Value = 4.32
$ws.Cells.Item(124, 21).Value = "05/11/2023 16:44"
$ws.Cells.Item(124, 22).Value = "https://www.betexplorer.com/football/hungary/merkantil-bank-liga/szeged-csanad-ga-budafoki-mte/dSqIBYEU/"

# New values for row 126 (previously held by row 124)
$ws.Cells.Item(126, 6).Value = "Csakvari"
$ws.Cells.Item(126, 8).Value = "Siofok"
$ws.Cells.Item(126, 10).Value = 2.04
$ws.Cells.Item(126, 12).Value = 2.02
$ws.Cells.Item(126, 13).Value = "05/11/2023 16:32"
$ws.Cells.Item(126, 14).Value = 3.2
$ws.Cells.Item(126, 16).Value = 3.49
$ws.Cells.Item(126, 17).Value = "05/11/2023 16:32"
$ws.Cells.Item(126, 18).Value = 3.36
$ws.Cells.Item(126, 20).Value = 3.62
$ws.Cells.Item(126, 21).Value = "05/11/2023 16:32"
$ws.Cells.Item(126, 22).Value = "https://www.betexplorer.com/football/hungary/merkantil-bank-liga/csakvari-siofok/0KYiJbGt/"

# --- Append the new Vasas vs Honved match as row 127.

# Seed the row first, then clone the formatting (bold/bordered index cell,
# date-formatted match-date cell) from row 126 so no new style entries are
# minted.
$ws.Cells.Item(127, 1).Value = 126
$ws.Cells.Item(127, 2).Value = "hungary"
$ws.Cells.Item(127, 3).Value = "merkantil-bank-liga"
$ws.Cells.Item(127, 4).Value = "2023-2024"
$ws.Cells.Item(127, 5).Value = 45236.83333333334
$ws.Cells.Item(127, 6).Value = "Vasas"
$ws.Cells.Item(127, 7).Value = 3
$ws.Cells.Item(127, 8).Value = "Honved"
$ws.Cells.Item(127, 9).Value = 0
$ws.Cells.Item(127, 10).Value = 1.75
$ws.Cells.Item(127, 11).Value = "05/11/2023 08:12"
$ws.Cells.Item(127, 12).Value = 1.79
$ws.Cells.Item(127, 13).Value = "06/11/2023 19:50"
$ws.Cells.Item(127, 14).Value = 3.58
$ws.Cells.Item(127, 15).Value = "05/11/2023 08:12"
$ws.Cells.Item(127, 16).Value = 3.74
$ws.Cells.Item(127, 17).Value = "06/11/2023 19:50"
$ws.Cells.Item(127, 18).Value = 4.13
$ws.Cells.Item(127, 19).Value = "05/11/2023 08:12"
$ws.Cells.Item(127, 20).Value = 4.34
$ws.Cells.Item(127, 21).Value = "06/11/2023 19:50"
$ws.Cells.Item(127, 22).Value = "https://www.betexplorer.com/football/hungary/merkantil-bank-liga/vasas-honved/Snh1FGF5/"

# Clone cell formatting for the styled columns (A: bold+border index style,
# E: date-time number format) from the row above so the workbook reuses the
# existing style records instead of minting new ones.
$ws.Cells.Item(126, 1).Copy()
$ws.Cells.Item(127, 1).PasteSpecial(-4122)
$ws.Cells.Item(126, 5).Copy()
$ws.Cells.Item(127, 5).PasteSpecial(-4122)

$excel.CutCopyMode = 0
